$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1772.5
$ws.Range("J112").Value = 2004.2858
$ws.Range("L112").Value = 6012.857400000001
$ws.Range("N112").Value = -8228.857400000001

$ws.Range("H129").Value = 858.8205
$ws.Range("I129").Value = 838
$ws.Range("J129").Value = 880.7368
$ws.Range("K129").Value = 2514
$ws.Range("L129").Value = 2642.2104
$ws.Range("M129").Value = 2486
$ws.Range("N129").Value = -12642.2104

$ws.Range("H137").Value = 7290.0234
$ws.Range("I137").Value = 10536.615
$ws.Range("K137").Value = 31609.845
$ws.Range("M137").Value = -29059.845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2088.6667
$ws.Range("I2").Value = 2088.6667
$ws.Range("K2").Value = 2088.6667
$ws.Range("M2").Value = -1975.6667

$ws.Range("H32").Value = 12512.946
$ws.Range("I32").Value = 12074.937
$ws.Range("J32").Value = 14800.333
$ws.Range("K32").Value = 12074.937
$ws.Range("L32").Value = 14800.333
$ws.Range("M32").Value = -11787.937
$ws.Range("N32").Value = -15374.333

$ws.Range("H45").Value = 58825950
$ws.Range("I45").Value = 111113224
$ws.Range("J45").Value = 2770.5
$ws.Range("K45").Value = 111113224
$ws.Range("L45").Value = 2770.5
$ws.Range("M45").Value = -111112847
$ws.Range("N45").Value = -3524.5

$ws.Range("H74").Value = 1533.5536
$ws.Range("I74").Value = 1292.5209
$ws.Range("J74").Value = 2979.75
$ws.Range("K74").Value = 1292.5209
$ws.Range("L74").Value = 2979.75
$ws.Range("M74").Value = -418.5209
$ws.Range("N74").Value = -4727.75

$ws.Range("H77").Value = 1533.5536
$ws.Range("I77").Value = 1292.5209
$ws.Range("J77").Value = 2979.75
$ws.Range("K77").Value = 6462.604499999999
$ws.Range("L77").Value = 14898.75
$ws.Range("M77").Value = -2094.604499999999
$ws.Range("N77").Value = -23634.75

$ws.Range("H113").Value = 49999.5
$ws.Range("J113").Value = 49999.5
$ws.Range("L113").Value = 49999.5
$ws.Range("N113").Value = -58677.5

$ws.Range("H116").Value = 2088.6667
$ws.Range("I116").Value = 2088.6667
$ws.Range("K116").Value = 2088.6667
$ws.Range("M116").Value = 205.3332999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2088.6667
$ws.Range("I3").Value = 2088.6667
$ws.Range("K3").Value = 2088.6667
$ws.Range("M3").Value = -1974.6667

$ws.Range("H94").Value = 1156.2778
$ws.Range("I94").Value = 1378.4546
$ws.Range("J94").Value = 807.1429000000001
$ws.Range("K94").Value = 1378.4546
$ws.Range("L94").Value = 807.1429000000001
$ws.Range("M94").Value = -927.4546
$ws.Range("N94").Value = -1709.1429

$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 519984.97
$ws.Range("I134").Value = 1299.9445
$ws.Range("J134").Value = 1557355
$ws.Range("K134").Value = 3899.8335
$ws.Range("L134").Value = 4672065
$ws.Range("M134").Value = -1364.8335
$ws.Range("N134").Value = -4677135

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 10944103
$ws.Range("I68").Value = 6536905
$ws.Range("J68").Value = 15626750
$ws.Range("K68").Value = 19610715
$ws.Range("L68").Value = 46880250
$ws.Range("M68").Value = -19609904
$ws.Range("N68").Value = -46881872

$ws.Range("H71").Value = 10944103
$ws.Range("I71").Value = 6536905
$ws.Range("J71").Value = 15626750
$ws.Range("K71").Value = 58832145
$ws.Range("L71").Value = 140640750
$ws.Range("M71").Value = -58828089
$ws.Range("N71").Value = -140648862

$ws.Range("H107").Value = 760.66
$ws.Range("I107").Value = 362.8889
$ws.Range("K107").Value = 1088.6667
$ws.Range("M107").Value = 831.3333

$ws.Range("H133").Value = 6169.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8697021
$ws.Range("I113").Value = 18183160
$ws.Range("J113").Value = 1393
$ws.Range("K113").Value = 18183160
$ws.Range("L113").Value = 1393
$ws.Range("M113").Value = -18180990
$ws.Range("N113").Value = -5733

$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226

$ws.Range("H126").Value = 9087.666999999999
$ws.Range("I126").Value = 36971
$ws.Range("J126").Value = 2116.8333
$ws.Range("K126").Value = 110913
$ws.Range("L126").Value = 6350.499899999999
$ws.Range("M126").Value = -108443
$ws.Range("N126").Value = -11290.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 58827092
$ws.Range("I7").Value = 100001976
$ws.Range("J7").Value = 5830
$ws.Range("K7").Value = 100001976
$ws.Range("L7").Value = 5830
$ws.Range("M7").Value = -100001864
$ws.Range("N7").Value = -6054

$ws.Range("H40").Value = 4452.467
$ws.Range("I40").Value = 3898.5
$ws.Range("J40").Value = 5560.4
$ws.Range("K40").Value = 3898.5
$ws.Range("L40").Value = 5560.4
$ws.Range("M40").Value = -3762.5
$ws.Range("N40").Value = -5832.4

$ws.Range("H61").Value = 1400
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 1400
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 1400
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -1804

$ws.Range("H113").Value = 1400
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1400
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5740

$ws.Range("H126").Value = 58827092
$ws.Range("I126").Value = 100001976
$ws.Range("J126").Value = 5830
$ws.Range("K126").Value = 300005928
$ws.Range("L126").Value = 17490
$ws.Range("M126").Value = -300003458
$ws.Range("N126").Value = -22430

$ws.Range("H136").Value = 1469.3793
$ws.Range("I136").Value = 1118.6666
$ws.Range("J136").Value = 2390
$ws.Range("K136").Value = 3355.9998
$ws.Range("L136").Value = 7170
$ws.Range("M136").Value = -805.9998000000001
$ws.Range("N136").Value = -12270

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 70011
$ws.Range("J29").Value = 70011
$ws.Range("L29").Value = 70011
$ws.Range("N29").Value = -70591

$ws.Range("H57").Value = 31500
$ws.Range("J57").Value = 31500
$ws.Range("L57").Value = 31500
$ws.Range("N57").Value = -33008

$ws.Range("H113").Value = 848.4
$ws.Range("I113").Value = 580.6667
$ws.Range("J113").Value = 1250
$ws.Range("K113").Value = 1742.0001
$ws.Range("L113").Value = 3750
$ws.Range("M113").Value = 427.9999
$ws.Range("N113").Value = -8090

$ws.Range("H126").Value = 2942702.8
$ws.Range("I126").Value = 3678002.8
$ws.Range("J126").Value = 1502.5
$ws.Range("K126").Value = 11034008.4
$ws.Range("L126").Value = 4507.5
$ws.Range("M126").Value = -11031538.4
$ws.Range("N126").Value = -9447.5

$ws.Range("H132").Value = 1972.15
$ws.Range("I132").Value = 1535.8
$ws.Range("J132").Value = 2699.4
$ws.Range("K132").Value = 4607.4
$ws.Range("L132").Value = 8098.200000000001
$ws.Range("M132").Value = -2077.4
$ws.Range("N132").Value = -13158.2
